$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-400) holds the "Förändrad" (changed) date as an Excel
# serial date. The update bumps that date by one day (45189 -> 45190,
# i.e. 2023-09-20 -> 2023-09-21) for every data row.
$ws.Range("C2:C400").Value = 45190
